$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Questionnaire": two new survey responses (rows 11 and 12),
# and clear the now-stray R10 answer.
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Questionnaire")

# Row 10 used to have a "likely to recommend" answer in R10; that answer
# was removed.
$ws.Range("R10").ClearContents()

# --- Row 11 (participant #9, "Oliver") ---
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 4
$ws.Range("D11").Value = "New experience of survey, never filled out a survey in this topic before, I would prefer more questions to get deeper."
$ws.Range("E11").Value = "–"
$ws.Range("F11").Value = "yes"
$ws.Range("G11").Value = "I think it helps reflecting the own relationship"
$ws.Range("H11").Value = "–"
$ws.Range("I11").Value = 6
$ws.Range("J11").Value = "The survey on the smartphone worked perfectly, but the form was confusing to fill out. "
$ws.Range("K11").Value = "–"
$ws.Range("L11").Value = 5
$ws.Range("M11").Value = "I liked the question, but it could have been more"
$ws.Range("N11").Value = "–"
$ws.Range("O11").Value = "no"
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = "For this, the survey should be larger, have more questions… And it would be good if there was the possibility of a rating at the end and then the option to compare with other participants (anonymous)."

# --- Row 12 (participant #10, "Hong") ---
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = ": I like the techn. background (design, etc) of the survey. What is missing is the option STUDENT as profession and no description if the survey is anonymous or not."
$ws.Range("E12").Value = "–"
$ws.Range("F12").Value = "yes"
$ws.Range("G12").Value = "maybe a little bit. It summarizes a relationship and motivate to improve it."
$ws.Range("H12").Value = "–"
$ws.Range("I12").Value = 6
$ws.Range("J12").Value = "design"
$ws.Range("K12").Value = "No, not enough questions available to help my relationship"
$ws.Range("L12").Value = 5
$ws.Range("M12").Value = "the idea of survey"
$ws.Range("N12").Value = "confusing feedback form (design)"
$ws.Range("O12").Value = "yes"
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = "I would pay for it, if the survey contains more personal questions. "

# Match the existing columnar formatting used by the rows above (B=id
# style, C/F/I/L/O=centered, P=currency format).
$ws.Range("B11").Style = $ws.Range("B10").Style
$ws.Range("B12").Style = $ws.Range("B10").Style
$ws.Range("C11:C12").Style = $ws.Range("C10").Style
$ws.Range("F11:F12").Style = $ws.Range("F10").Style
$ws.Range("I11:I12").Style = $ws.Range("I10").Style
$ws.Range("L11:L12").Style = $ws.Range("L10").Style
$ws.Range("O11:O12").Style = $ws.Range("O10").Style
$ws.Range("P11:P12").Style = $ws.Range("P10").Style

# The free-text answers that were pasted in straight from a GitHub
# comment box came in with GitHub's own font/colour; reproduce that
# look (Helvetica 12pt, #24292E) and carry it onto the rest of the
# pasted-in answers.
$font = $ws.Range("D11").Font
$font.Name = "Helvetica"
$font.Color = 3025188

$ws.Range("D11").Copy() | Out-Null
$ws.Range("G11,Q11,D12,G12,K12,Q12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ----------------------------------------------------------------------
# Sheet "Participants": the two new respondents.
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Participants")
$ws2.Range("B11").Value = 9
$ws2.Range("C11").Value = "Oliver"
$ws2.Range("B12").Value = 10
$ws2.Range("C12").Value = "Hong"
